# Applies the "Forgot to save before commit" edit:
#  - switches the active tab from "Partants" to "All"
#  - updates the "All" sheet's field-mapping table (column C / some column B
#    labels) to add the newly-scraped fields (odds, commentary, race info,
#    pedigree, earnings split by year, etc.)
#  - updates the selection shown on the "All" sheet

$wb = $excel.ActiveWorkbook
$wsAll = $wb.Worksheets.Item("All")

# --- Column C (internal field name) / column B (French label) updates ---

# "Gains" row loses its internal-name mapping (was "earnings", now handled
# elsewhere as earnings_*).
$wsAll.Range("C11").ClearContents()

# "Performances" row's internal name changes from "places" to "history".
$wsAll.Range("C13").Value = "history"

# "Rapports prob. Last" label becomes "Rapports prob. Direct".
$wsAll.Range("B15").Value = "Rapports prob. Direct"

# Reference odds row gains its internal field name.
$wsAll.Range("C30").Value = "odds_ref"

# "Rapports prob. Last" (Arrivée block) becomes "Rapports prob. Direct" and
# gains its internal field name.
$wsAll.Range("B31").Value = "Rapports prob. Direct"
$wsAll.Range("C31").Value = "odds_direct"

# Post-race commentary row gains its internal field name.
$wsAll.Range("C32").Value = "commentary"

# "Deep" (horse record) block: race/pedigree fields gain internal names.
$wsAll.Range("C39").Value = "race"
$wsAll.Range("C42").Value = "coat"
$wsAll.Range("C43").Value = "races_run"
$wsAll.Range("C44").Value = "victories"
$wsAll.Range("C45").Value = "places"
$wsAll.Range("C46").Value = "earnings_career"
$wsAll.Range("C47").Value = "earnings_last_year"
$wsAll.Range("C48").Value = "earnings_victory"

# "Gains N-1" (4th occurrence) becomes "Gains N" and gains an internal name.
$wsAll.Range("B49").Value = "Gains N"
$wsAll.Range("C49").Value = "earnings_current_year"

# Owner/breeder/pedigree fields gain internal names.
$wsAll.Range("C51").Value = "owner"
$wsAll.Range("C52").Value = "breeder"
$wsAll.Range("C53").Value = "father"
$wsAll.Range("C54").Value = "mother"
$wsAll.Range("C55").Value = "mother's father"

# --- View-state: make "All" the active/visible tab with its new selection ---

$wsAll.Activate()
$wsAll.Range("C3:C22").Select()
